$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4932.905
$ws.Range("I33").Value = 91.61539
$ws.Range("J33").Value = 12800
$ws.Range("K33").Value = 91.61539
$ws.Range("L33").Value = 12800
$ws.Range("M33").Value = 137.38461
$ws.Range("N33").Value = -13258
$ws.Range("H62").Value = 2360.087
$ws.Range("I62").Value = 2188.889
$ws.Range("J62").Value = 2976.4
$ws.Range("K62").Value = 2188.889
$ws.Range("L62").Value = 2976.4
$ws.Range("M62").Value = -1564.889
$ws.Range("N62").Value = -4224.4
$ws.Range("H65").Value = 2360.087
$ws.Range("I65").Value = 2188.889
$ws.Range("J65").Value = 2976.4
$ws.Range("K65").Value = 10944.445
$ws.Range("L65").Value = 14882
$ws.Range("M65").Value = -7824.445
$ws.Range("N65").Value = -21122
$ws.Range("H98").Value = 3476019.2
$ws.Range("I98").Value = 4421.3076
$ws.Range("J98").Value = 18519610
$ws.Range("K98").Value = 4421.3076
$ws.Range("L98").Value = 18519610
$ws.Range("M98").Value = -2923.3076
$ws.Range("N98").Value = -18522606
$ws.Range("H106").Value = 2917.5
$ws.Range("I106").Value = 2917.5
$ws.Range("K106").Value = 2917.5
$ws.Range("M106").Value = -2286.5
$ws.Range("H122").Value = 3476019.2
$ws.Range("I122").Value = 4421.3076
$ws.Range("J122").Value = 18519610
$ws.Range("K122").Value = 13263.9228
$ws.Range("L122").Value = 55558830
$ws.Range("M122").Value = -10813.9228
$ws.Range("N122").Value = -55563730
$ws.Range("H138").Value = 4135.377
$ws.Range("I138").Value = 3433.9285
$ws.Range("J138").Value = 4344.3193
$ws.Range("K138").Value = 10301.7855
$ws.Range("L138").Value = 13032.9579
$ws.Range("M138").Value = -5161.7855
$ws.Range("N138").Value = -23312.9579

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1909.579
$ws.Range("I2").Value = 1182.72
$ws.Range("J2").Value = 3307.3845
$ws.Range("K2").Value = 1182.72
$ws.Range("L2").Value = 3307.3845
$ws.Range("M2").Value = -1069.72
$ws.Range("N2").Value = -3533.3845
$ws.Range("H31").Value = 4034.1538
$ws.Range("I31").Value = 1493.75
$ws.Range("K31").Value = 1493.75
$ws.Range("M31").Value = -1199.75
$ws.Range("H32").Value = 862513.9
$ws.Range("I32").Value = 9443.276
$ws.Range("J32").Value = 10886094
$ws.Range("K32").Value = 9443.276
$ws.Range("L32").Value = 10886094
$ws.Range("M32").Value = -9156.276
$ws.Range("N32").Value = -10886668
$ws.Range("H61").Value = 2533.4783
$ws.Range("I61").Value = 2582.75
$ws.Range("J61").Value = 2205
$ws.Range("K61").Value = 2582.75
$ws.Range("L61").Value = 2205
$ws.Range("M61").Value = -2370.75
$ws.Range("N61").Value = -2629
$ws.Range("H63").Value = 4085.2942
$ws.Range("I63").Value = 1431.25
$ws.Range("J63").Value = 6444.4443
$ws.Range("K63").Value = 1431.25
$ws.Range("L63").Value = 6444.4443
$ws.Range("M63").Value = -745.25
$ws.Range("N63").Value = -7816.4443
$ws.Range("H66").Value = 4085.2942
$ws.Range("I66").Value = 1431.25
$ws.Range("J66").Value = 6444.4443
$ws.Range("K66").Value = 7156.25
$ws.Range("L66").Value = 32222.2215
$ws.Range("M66").Value = -3724.25
$ws.Range("N66").Value = -39086.2215
$ws.Range("H116").Value = 1909.579
$ws.Range("I116").Value = 1182.72
$ws.Range("J116").Value = 3307.3845
$ws.Range("K116").Value = 1182.72
$ws.Range("L116").Value = 3307.3845
$ws.Range("M116").Value = 1111.28
$ws.Range("N116").Value = -7895.3845
$ws.Range("H122").Value = 22095.527
$ws.Range("I122").Value = 25067.227
$ws.Range("J122").Value = 3671
$ws.Range("K122").Value = 75201.681
$ws.Range("L122").Value = 11013
$ws.Range("M122").Value = -72751.681
$ws.Range("N122").Value = -15913
$ws.Range("H136").Value = 2533.4783
$ws.Range("I136").Value = 2582.75
$ws.Range("J136").Value = 2205
$ws.Range("K136").Value = 7748.25
$ws.Range("L136").Value = 6615
$ws.Range("M136").Value = -5198.25
$ws.Range("N136").Value = -11715

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1909.579
$ws.Range("I3").Value = 1182.72
$ws.Range("J3").Value = 3307.3845
$ws.Range("K3").Value = 1182.72
$ws.Range("L3").Value = 3307.3845
$ws.Range("M3").Value = -1068.72
$ws.Range("N3").Value = -3535.3845
$ws.Range("H99").Value = 1621.6316
$ws.Range("I99").Value = 1109.091
$ws.Range("J99").Value = 2326.375
$ws.Range("K99").Value = 1109.091
$ws.Range("L99").Value = 2326.375
$ws.Range("M99").Value = 388.9090000000001
$ws.Range("N99").Value = -5322.375
$ws.Range("H102").Value = 12358.857
$ws.Range("I102").Value = 12358.857
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 12358.857
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -9113.857
$ws.Range("N102").ClearContents()
$ws.Range("H134").Value = 7361.1816
$ws.Range("I134").Value = 746.4167
$ws.Range("J134").Value = 15298.9
$ws.Range("K134").Value = 2239.2501
$ws.Range("L134").Value = 45896.7
$ws.Range("M134").Value = 295.7498999999998
$ws.Range("N134").Value = -50966.7

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H58").Value = 1202.1428
$ws.Range("I58").Value = 703.0625
$ws.Range("J58").Value = 1867.5834
$ws.Range("K58").Value = 703.0625
$ws.Range("L58").Value = 1867.5834
$ws.Range("M58").Value = -500.0625
$ws.Range("N58").Value = -2273.5834
$ws.Range("H62").Value = 4053.9565
$ws.Range("I62").Value = 3975.9092
$ws.Range("K62").Value = 3975.9092
$ws.Range("M62").Value = -3351.9092
$ws.Range("H65").Value = 4053.9565
$ws.Range("I65").Value = 3975.9092
$ws.Range("K65").Value = 19879.546
$ws.Range("M65").Value = -16759.546
$ws.Range("H68").Value = 15937.833
$ws.Range("J68").Value = 15937.833
$ws.Range("L68").Value = 15937.833
$ws.Range("N68").Value = -17435.833
$ws.Range("H71").Value = 15937.833
$ws.Range("J71").Value = 15937.833
$ws.Range("L71").Value = 47813.499
$ws.Range("N71").Value = -55301.499
$ws.Range("H105").Value = 1201.3077
$ws.Range("I105").Value = 1065.2
$ws.Range("J105").Value = 1655
$ws.Range("K105").Value = 1065.2
$ws.Range("L105").Value = 1655
$ws.Range("M105").Value = 681.8
$ws.Range("N105").Value = -5149
$ws.Range("H132").Value = 2762.5386
$ws.Range("I132").Value = 1500.1666
$ws.Range("J132").Value = 3844.5715
$ws.Range("K132").Value = 4500.4998
$ws.Range("L132").Value = 11533.7145
$ws.Range("M132").Value = -1970.4998
$ws.Range("N132").Value = -16593.7145
$ws.Range("H134").Value = 674.7778
$ws.Range("I134").Value = 685.25
$ws.Range("J134").Value = 644.8570999999999
$ws.Range("K134").Value = 2055.75
$ws.Range("L134").Value = 1934.5713
$ws.Range("M134").Value = 479.25
$ws.Range("N134").Value = -7004.5713
$ws.Range("H136").Value = 1202.1428
$ws.Range("I136").Value = 703.0625
$ws.Range("J136").Value = 1867.5834
$ws.Range("K136").Value = 2109.1875
$ws.Range("L136").Value = 5602.7502
$ws.Range("M136").Value = 440.8125
$ws.Range("N136").Value = -10702.7502

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2294
$ws.Range("I55").Value = 704
$ws.Range("J55").Value = 2400
$ws.Range("K55").Value = 2112
$ws.Range("L55").Value = 7200
$ws.Range("M55").Value = -1935
$ws.Range("N55").Value = -7554
$ws.Range("H68").Value = 951.3
$ws.Range("I68").Value = 800.25
$ws.Range("J68").Value = 1052
$ws.Range("K68").Value = 2400.75
$ws.Range("L68").Value = 3156
$ws.Range("M68").Value = -1589.75
$ws.Range("N68").Value = -4778
$ws.Range("H71").Value = 951.3
$ws.Range("I71").Value = 800.25
$ws.Range("J71").Value = 1052
$ws.Range("K71").Value = 7202.25
$ws.Range("L71").Value = 9468
$ws.Range("M71").Value = -3146.25
$ws.Range("N71").Value = -17580
$ws.Range("H131").Value = 15152823
$ws.Range("J131").Value = 16667742
$ws.Range("L131").Value = 50003226
$ws.Range("N131").Value = -50013306

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3258.2856
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 3561.6
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 10684.8
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -15584.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 22225234
$ws.Range("I40").Value = 2996
$ws.Range("K40").Value = 2996
$ws.Range("M40").Value = -2860
$ws.Range("H122").Value = 4064.2
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4160.5
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 12481.5
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -17381.5

Write-Host "Applied all Durandal_Profits market price updates."